# Final "Thank you" slide: the title placeholder currently holds an empty
# paragraph (just an endParaRPr). Fill it in with the closing "Thank you"
# text, as a normal PowerPoint user would by clicking into the title and
# typing.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item($p.Slides.Count)
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Thank you"
